# liensLogistiqueFr.xlsx — remove the four "Liens Ressources Humaines :" (HR
# links) mini-blocks (K'IAM / SMARTRH rows) that appear under each
# "Logistique ... ClientFR" section on sheet Feuil1.
#
# Each block occupies two rows, columns A:D:
#   row N   : A=<category>, B="Liens Ressources Humaines :", C="K'IAM",    D=<K'IAM url>
#   row N+1 : A=<category>, B=<blank>,                        C="SMARTRH", D=<SMARTRH url>
# The blocks start at rows 15, 41, 67 and 91.
#
# The cells are cleared (contents only — formatting/style stays untouched,
# matching the <c r="..." s=".."/> cells left behind in the saved file) and
# the hyperlink attached to each K'IAM cell (column D of the block's first
# row) is removed along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blockStartRows = @(15, 41, 67, 91)

foreach ($startRow in $blockStartRows) {
    # Drop the hyperlink that lives on the K'IAM row (column D) of this
    # block before the text backing it disappears.
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Row -eq $startRow) {
            $hl.Delete()
        }
    }

    $endRow = $startRow + 1
    $ws.Range("A" + $startRow + ":D" + $endRow).ClearContents()
}

# Best-effort cosmetic bookkeeping to mirror the author's final cursor
# position (does not affect any cell content/values).
$ws.Range("D92").Select()
try {
    $excel.ActiveWindow.ScrollRow = 40
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}

$wb.Save()
